$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells we touch are written as text so that
# numeric-looking strings keep their exact original formatting/precision.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D40","D41","D42","D44","D45","D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Rows 9-28: coin list shifted up by one position (WazirX moved down the list) ---
$ws.Range("B9").Value = "BitrueCoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D9").Value = "0.02881"
$ws.Range("E9").Value = "8BitrueCoinBTR"

$ws.Range("B10").Value = "BitMartToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D10").Value = "0.09411"
$ws.Range("E10").Value = "9BitMartTokenBMX"

$ws.Range("B11").Value = "BitForexToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D11").Value = "0.001521"
$ws.Range("E11").Value = "10BitForexTokenBF"

$ws.Range("B12").Value = "One"
$ws.Range("C12").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D12").Value = "0.0006000"
$ws.Range("E12").Value = "11OneONE"

$ws.Range("B13").Value = "TigerCash"
$ws.Range("C13").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D13").Value = "0.006187"
$ws.Range("E13").Value = "12TigerCashTCH"

$ws.Range("B14").Value = "LEO"
$ws.Range("C14").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D14").Value = "3.589"
$ws.Range("E14").Value = "13LEOLEO"

$ws.Range("B15").Value = "GateToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D15").Value = "3.016"
$ws.Range("E15").Value = "14GateTokenGT"

$ws.Range("B16").Value = "BTSEToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D16").Value = "2.118"
$ws.Range("E16").Value = "15BTSETokenBTSE"

$ws.Range("B17").Value = "BitpandaEcosystemToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D17").Value = "0.3157"
$ws.Range("E17").Value = "16BitpandaEcosystemTokenBEST"

$ws.Range("B18").Value = "WazirX"
$ws.Range("C18").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D18").Value = "0.1346"
$ws.Range("E18").Value = "17WazirXWRX"

$ws.Range("B19").Value = "MandalaExchangeToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D19").Value = "0.07001"
$ws.Range("E19").Value = "18MandalaExchangeTokenMDX"

$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "0.03184"
$ws.Range("E20").Value = "19LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1320"
$ws.Range("E21").Value = "20ProBitTokenPROB"

$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "3.745"
$ws.Range("E22").Value = "21MCDexMCB"

$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "0.04661"
$ws.Range("E23").Value = "22CoinExTokenCET"

$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").Value = "0.1350"
$ws.Range("E24").Value = "23ZBTokenZB"

$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").Value = "0.001252"
$ws.Range("E25").Value = "24BitKanKAN"

$ws.Range("B26").Value = "HotbitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D26").Value = "0.004601"
$ws.Range("E26").Value = "25HotbitTokenHTB"

$ws.Range("B27").Value = "NitroEx"
$ws.Range("C27").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D27").Value = "0.00009601"
$ws.Range("E27").Value = "26NitroExNTX"

$ws.Range("B28").Value = "UpBots"
$ws.Range("C28").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D28").Value = "0.0001938"
$ws.Range("E28").Value = "27UpBotsUBXT"

# --- Price-only refreshes ---
$ws.Range("D2").Value = "246.98"
$ws.Range("D3").Value = "26.37"
$ws.Range("D4").Value = "5.083"
$ws.Range("D5").Value = "0.05615"
$ws.Range("D6").Value = "6.515"
$ws.Range("D7").Value = "0.8132"
$ws.Range("D8").Value = "0.8456"
$ws.Range("D40").Value = "0.03679"
$ws.Range("D42").Value = "0.1059"
$ws.Range("D44").Value = "0.008916"
$ws.Range("D45").Value = "0.00005296"
$ws.Range("D48").Value = "0.002528"

# --- Row 41: price + label refresh (Best in 24h tag) ---
$ws.Range("D41").Value = "0.006141"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

# --- Row 47: label refresh (Worst in 24h tag) ---
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

